# Fixed limit switch pin labels
#   "Limit switch 2"            -> "Top limit switch"    (table 3, row 10, col 4)
#   "Limit switch " + "1"       -> "Bottom limit switch"  (table 3, row 11, col 4)
# Both target cells end up holding the new text split across two runs,
# e.g. "Top l" / "imit switch", matching the authored edit.

$d = $word.ActiveDocument

function Set-CellSplitText($Table, $Row, $Col, $FirstPart, $SecondPart) {
    $newText = $FirstPart + $SecondPart

    # 1) Overwrite the cell's paragraph text in one shot. Word COM only
    #    rewrites the text of the first run a range overlaps, so if the
    #    cell originally held more than one run, any trailing run(s) keep
    #    their old characters tacked on after our new text - clean that up
    #    below.
    $cell = $Table.Cell($Row, $Col)
    $full = $cell.Range
    $start = $full.Start
    $end = $full.End

    $target = $d.Range($start, $end - 1)
    $target.Text = $newText

    $cell = $Table.Cell($Row, $Col)
    $full = $cell.Range
    $expectedEnd = $start + $newText.Length + 1
    if ($full.End -gt $expectedEnd) {
        $leftover = $d.Range($start + $newText.Length, $full.End - 1)
        $leftover.Text = ""
    }

    # 2) Split the now-single run into two runs at the boundary between
    #    FirstPart and SecondPart, matching the authored commit. Toggling a
    #    character-level property (and toggling it straight back) on just
    #    the second part forces the engine to materialize it as its own run.
    $cell = $Table.Cell($Row, $Col)
    $full = $cell.Range
    $splitAt = $start + $FirstPart.Length
    $secondRun = $d.Range($splitAt, $full.End - 1)
    $secondRun.Bold = 1
    $secondRun.Bold = 0
}

$pinTable = $d.Tables.Item(3)

Set-CellSplitText $pinTable 10 4 "Top l" "imit switch"
Set-CellSplitText $pinTable 11 4 "Bottom l" "imit switch"
